$d = $word.ActiveDocument

# Update the date line
$d.Paragraphs.Item(1).Range.Find.Execute("2025-09-13 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-14 Sunday", 2) | Out-Null

# Update table cells (20 rows x 5 columns), addressed positionally to avoid
# ambiguity from duplicate old values (e.g. "9-1=8" appears twice).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "30+45=75"
$t.Cell(1,2).Range.Text = "42-39=3"
$t.Cell(1,3).Range.Text = "34-6=28"
$t.Cell(1,4).Range.Text = "94-57=37"
$t.Cell(1,5).Range.Text = "68-24=44"
$t.Cell(2,1).Range.Text = "44+0=44"
$t.Cell(2,2).Range.Text = "11+0=11"
$t.Cell(2,3).Range.Text = "93-39=54"
$t.Cell(2,4).Range.Text = "63-30=33"
$t.Cell(2,5).Range.Text = "23+67=90"
$t.Cell(3,1).Range.Text = "80-50=30"
$t.Cell(3,2).Range.Text = "61-44=17"
$t.Cell(3,3).Range.Text = "62-18=44"
$t.Cell(3,4).Range.Text = "73-5=68"
$t.Cell(3,5).Range.Text = "31+31=62"
$t.Cell(4,1).Range.Text = "33+11=44"
$t.Cell(4,2).Range.Text = "73-63=10"
$t.Cell(4,3).Range.Text = "35+41=76"
$t.Cell(4,4).Range.Text = "40-9=31"
$t.Cell(4,5).Range.Text = "31-4=27"
$t.Cell(5,1).Range.Text = "53-10=43"
$t.Cell(5,2).Range.Text = "93-27=66"
$t.Cell(5,3).Range.Text = "91-28=63"
$t.Cell(5,4).Range.Text = "96-93=3"
$t.Cell(5,5).Range.Text = "4+8=12"
$t.Cell(6,1).Range.Text = "1+65=66"
$t.Cell(6,2).Range.Text = "37+25=62"
$t.Cell(6,3).Range.Text = "15-9=6"
$t.Cell(6,4).Range.Text = "31+41=72"
$t.Cell(6,5).Range.Text = "55-27=28"
$t.Cell(7,1).Range.Text = "61-35=26"
$t.Cell(7,2).Range.Text = "62+13=75"
$t.Cell(7,3).Range.Text = "38+30=68"
$t.Cell(7,4).Range.Text = "71+26=97"
$t.Cell(7,5).Range.Text = "12+58=70"
$t.Cell(8,1).Range.Text = "66+24=90"
$t.Cell(8,2).Range.Text = "48+40=88"
$t.Cell(8,3).Range.Text = "93-44=49"
$t.Cell(8,4).Range.Text = "69-49=20"
$t.Cell(8,5).Range.Text = "60-56=4"
$t.Cell(9,1).Range.Text = "31-4=27"
$t.Cell(9,2).Range.Text = "30+12=42"
$t.Cell(9,3).Range.Text = "36+35=71"
$t.Cell(9,4).Range.Text = "5+24=29"
$t.Cell(9,5).Range.Text = "16+1=17"
$t.Cell(10,1).Range.Text = "67-28=39"
$t.Cell(10,2).Range.Text = "64+20=84"
$t.Cell(10,3).Range.Text = "53-9=44"
$t.Cell(10,4).Range.Text = "55+5=60"
$t.Cell(10,5).Range.Text = "90-84=6"
$t.Cell(11,1).Range.Text = "7+51=58"
$t.Cell(11,2).Range.Text = "27-9=18"
$t.Cell(11,3).Range.Text = "63+35=98"
$t.Cell(11,4).Range.Text = "6+33=39"
$t.Cell(11,5).Range.Text = "24-3=21"
$t.Cell(12,1).Range.Text = "4+73=77"
$t.Cell(12,2).Range.Text = "91-48=43"
$t.Cell(12,3).Range.Text = "74-70=4"
$t.Cell(12,4).Range.Text = "88-84=4"
$t.Cell(12,5).Range.Text = "16+48=64"
$t.Cell(13,1).Range.Text = "45-35=10"
$t.Cell(13,2).Range.Text = "28+60=88"
$t.Cell(13,3).Range.Text = "60+12=72"
$t.Cell(13,4).Range.Text = "50-16=34"
$t.Cell(13,5).Range.Text = "53+4=57"
$t.Cell(14,1).Range.Text = "12+22=34"
$t.Cell(14,2).Range.Text = "91-4=87"
$t.Cell(14,3).Range.Text = "93-89=4"
$t.Cell(14,4).Range.Text = "6+61=67"
$t.Cell(14,5).Range.Text = "82-33=49"
$t.Cell(15,1).Range.Text = "91-22=69"
$t.Cell(15,2).Range.Text = "51+21=72"
$t.Cell(15,3).Range.Text = "75-32=43"
$t.Cell(15,4).Range.Text = "66+19=85"
$t.Cell(15,5).Range.Text = "85-9=76"
$t.Cell(16,1).Range.Text = "61-14=47"
$t.Cell(16,2).Range.Text = "59-31=28"
$t.Cell(16,3).Range.Text = "7+14=21"
$t.Cell(16,4).Range.Text = "68-32=36"
$t.Cell(16,5).Range.Text = "21+5=26"
$t.Cell(17,1).Range.Text = "85+6=91"
$t.Cell(17,2).Range.Text = "59-58=1"
$t.Cell(17,3).Range.Text = "75-67=8"
$t.Cell(17,4).Range.Text = "49+22=71"
$t.Cell(17,5).Range.Text = "62+5=67"
$t.Cell(18,1).Range.Text = "77-0=77"
$t.Cell(18,2).Range.Text = "32-4=28"
$t.Cell(18,3).Range.Text = "64-35=29"
$t.Cell(18,4).Range.Text = "45-14=31"
$t.Cell(18,5).Range.Text = "18+0=18"
$t.Cell(19,1).Range.Text = "66+11=77"
$t.Cell(19,2).Range.Text = "68-66=2"
$t.Cell(19,3).Range.Text = "63-15=48"
$t.Cell(19,4).Range.Text = "47-38=9"
$t.Cell(19,5).Range.Text = "46-44=2"
$t.Cell(20,1).Range.Text = "67-20=47"
$t.Cell(20,2).Range.Text = "60-22=38"
$t.Cell(20,3).Range.Text = "24-18=6"
$t.Cell(20,4).Range.Text = "10+10=20"
$t.Cell(20,5).Range.Text = "49-0=49"
